$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "1 Us-en" -> "Add_Labels"
$ws.Name = "Add_Labels"

# Update the header row to the new label columns
$ws.Range("A1").Value = "Label_name"
$ws.Range("B1").Value = "1 US-en"
$ws.Range("C1").Value = "2 IN-hi"

# Remove the old demo-label data rows (rows 2 and 3) entirely so the
# used range shrinks back down to just the header row
$ws.Range("A2:C3").EntireRow.Delete()

# Swap the column widths: column A grows to 28, columns B:C shrink to 20
# (ColumnWidth adds the standard ~0.8333 char padding on write, so back it
# out here to land on the exact target stored width)
$ws.Columns.Item(1).ColumnWidth = 28 - 5/6
$ws.Range("B:C").ColumnWidth = 20 - 5/6
